$d = $word.ActiveDocument

$replacements = @(
    @("266÷3=", "921÷9="),
    @("220÷5=", "212÷2="),
    @("483÷4=", "424÷6="),
    @("940÷9=", "695÷7="),
    @("923÷4=", "793÷9="),
    @("560÷7=", "621÷8="),
    @("256÷6=", "399÷2="),
    @("158÷2=", "787÷4="),
    @("847÷3=", "433÷6="),
    @("844÷6=", "971÷6="),
    @("936÷6=", "961÷2="),
    @("695÷9=", "291÷3="),
    @("670÷2=", "290÷6="),
    @("536÷3=", "568÷2="),
    @("461÷6=", "265÷2="),
    @("414÷6=", "963÷3="),
    @("101÷5=", "136÷4="),
    @("533÷9=", "684÷8="),
    @("750÷3=", "284÷3="),
    @("765÷6=", "155÷8="),
    @("626÷7=", "978÷7="),
    @("519÷3=", "175÷2="),
    @("701÷3=", "456÷3="),
    @("687÷5=", "427÷2="),
    @("464÷6=", "505÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
